# Apply checklist updates to the "Checkliste Dokumente" workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 13 previously asked about empty filler pages ("Leere Füllseiten im
# Dokument?" / "nein"). It now asks whether all pages are filled with
# meaningful content, and the answer is "ja".
$ws.Range("B13").Value = "Alle Seiten mit sinnvollem Inhalt gefüllt?"
$ws.Range("D13").Value = "ja"

# Update the active selection on the sheet to C20.
$ws.Range("C20").Select()
